# Apply "Modif oublié de push cost" changes to the ST BOM sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / value corrections in existing rows ---------------------------

# Row 4 - Spacer comment shortened
$ws.Range("E4").Value = "Between steering wheel and Quick release"

# Row 5 - trailing period removed
$ws.Range("E5").Value = "Removing part of the Quick release"

# Row 7 - part renamed, trailing period removed
$ws.Range("C7").Value = "Quick Release Shaft"
$ws.Range("E7").Value = "Fixed part of the Quick release"

# Row 8 - part renamed, trailing period removed
$ws.Range("C8").Value = "Steering Shaft Pivot"
$ws.Range("E8").Value = "Bearing seat for the steering pivot"

# Row 9 - comment reworded
$ws.Range("E9").Value = "U-joint for steering column and ST_02002"

# Row 11 - comment reworded
$ws.Range("E11").Value = "Connect the steering column to the rack"

# Row 12 - part renamed, trailing period removed
$ws.Range("C12").Value = "Bearing, Ball, Radial"
$ws.Range("E12").Value = "Steering pivot bearings"

# Row 17 - tie rod tapped insert split into right hand version
$ws.Range("C17").Value = "Tapped insert, right hand"
$ws.Range("E17").Value = "Right-hand thread, glued to carbon tube"
$ws.Range("F17").Value = 2

# Row 18 - left hand tapped insert (was a leftover template row)
$ws.Range("C18").Value = "Tapped insert, left hand"
$ws.Range("D18").Value = "m"
$ws.Range("E18").Value = "Left-hand thread, glued to carbon tube"
$ws.Range("F18").Value = 2

# Row 19 - carbon tube (was a leftover template row)
$ws.Range("C19").Value = "Carbon tube"
$ws.Range("D19").Value = "b"
$ws.Range("E19").Value = "carbon tubes for tie rod"
$ws.Range("F19").Value = 2

# Row 20 - rod end bearing, right hand (was a leftover template row)
$ws.Range("C20").Value = "Rod ends bearing, male r"
$ws.Range("D20").Value = "b"
$ws.Range("E20").Value = "2 with a left-hand thread"
$ws.Range("F20").Value = 2

# Row 21 - rod end bearing, left hand (was a leftover template row)
$ws.Range("C21").Value = "Rod ends bearing, male l"
$ws.Range("D21").Value = "b"
$ws.Range("E21").Value = "2 with a right-hand thread"
$ws.Range("F21").Value = 2

# Row 22 - M6 16mm spacer (was a leftover template row)
$ws.Range("C22").Value = "Spacer "
$ws.Range("D22").Value = "m"
$ws.Range("E22").Value = "M6 type 16 mm spacer"
$ws.Range("F22").Value = 2

# --- New row 23: M6 25mm spacer ------------------------------------------
$ws.Range("A22:G22").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 15

$ws.Range("C23").Value = "Spacer "
$ws.Range("D23").Value = "m"
$ws.Range("E23").Value = "M6 type 25 mm spacer"
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = "ST_04007"

# --- Row height fix-up (rows whose wrapped text now fits one line) -------
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15

# --- Column width tweaks ---------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 23.830729166666668
$ws.Columns.Item(5).ColumnWidth = 40.83072916666667

# --- Restore the selected cell as last used by the author -----------------
[void]$ws.Range("H20").Select()
